$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, H(category), I(cond_cat), K(correct_answer), L(stimulus),
#             M(conceptual), N(perceptual), O(typicality), P(n),
#             Q(p_typicality), R(p_conceptual), S(p_perceptual),
#             T(r_typicality), U(r_conceptual), V(r_perceptual)
$data = @(
    @(2, 'living_rooms', 'target', 'j', 'stimuli/img_5jp4f.png', 84.85714285714286, 67.83333333333333, 76.3452380952381, 42, 9, 9, 9, 8, 8, 9),
    @(3, 'living_rooms', 'target', 'j', 'stimuli/img_6a0hu.png', 61.275, 42.025, 51.65, 40, 4, 4, 4, 5, 4, 5),
    @(4, 'bedrooms', 'distractor', 'f', 'stimuli/img_md6k3.png', 81.57142857142857, 63.66666666666666, 72.61904761904762, 42, 8, 8, 8, 8, 8, 8),
    @(5, 'living_rooms', 'target', 'j', 'stimuli/img_j4ttn.png', 12.61904761904762, 11.42857142857143, 12.02380952380952, 42, 1, 1, 1, 1, 1, 1),
    @(6, 'bedrooms', 'distractor', 'f', 'stimuli/img_ymgcb.png', 83.73684210526316, 61.13157894736842, 72.43421052631579, 38, 8, 8, 8, 8, 8, 8),
    @(7, 'bedrooms', 'distractor', 'f', 'stimuli/img_swcci.png', 49.82926829268293, 28.46341463414634, 39.14634146341464, 41, 2, 2, 2, 2, 2, 2),
    @(8, 'living_rooms', 'target', 'j', 'stimuli/img_16kib.png', 80.97727272727273, 61.11363636363637, 71.04545454545455, 44, 8, 8, 8, 7, 7, 7),
    @(9, 'kitchens', 'distractor', 'f', 'stimuli/img_gztbt.png', 55.06451612903226, 26.09677419354839, 40.58064516129032, 31, 2, 2, 2, 2, 2, 2),
    @(10, 'kitchens', 'distractor', 'f', 'stimuli/img_xtb79.png', 55.59375, 31.40625, 43.5, 32, 2, 2, 2, 2, 2, 2),
    @(12, 'bedrooms', 'distractor', 'f', 'stimuli/img_7caxh.png', 83.38095238095238, 64.26190476190476, 73.82142857142857, 42, 8, 8, 8, 8, 8, 8),
    @(13, 'bedrooms', 'distractor', 'f', 'stimuli/img_gv750.png', 83.51428571428572, 61.88571428571429, 72.7, 35, 8, 8, 8, 8, 8, 8),
    @(14, 'bedrooms', 'distractor', 'f', 'stimuli/img_b2jli.png', 83.15625, 63.8125, 73.484375, 32, 8, 8, 8, 8, 8, 8),
    @(15, 'bedrooms', 'distractor', 'f', 'stimuli/img_awuev.png', 44.21052631578947, 24.26315789473684, 34.23684210526316, 38, 2, 2, 2, 2, 2, 2),
    @(16, 'living_rooms', 'target', 'j', 'stimuli/img_z4jxm.png', 88.30952380952381, 72.64285714285714, 80.47619047619048, 42, 10, 10, 10, 10, 10, 10),
    @(17, 'living_rooms', 'target', 'j', 'stimuli/img_pdzf1.png', 86.23913043478261, 67.17391304347827, 76.70652173913044, 46, 9, 9, 9, 9, 9, 8),
    @(18, 'living_rooms', 'target', 'j', 'stimuli/img_ys3qz.png', 46.79545454545455, 31.20454545454545, 39, 44, 2, 2, 2, 2, 2, 2),
    @(19, 'kitchens', 'distractor', 'f', 'stimuli/img_463mq.png', 51.35294117647059, 30.20588235294118, 40.77941176470588, 34, 2, 2, 2, 2, 2, 2),
    @(20, 'bedrooms', 'distractor', 'f', 'stimuli/img_91csq.png', 50.44736842105263, 28.34210526315789, 39.39473684210526, 38, 2, 2, 2, 2, 2, 2),
    @(21, 'living_rooms', 'target', 'j', 'stimuli/img_vgh2g.png', 93.81395348837209, 78.27906976744185, 86.04651162790697, 43, 10, 10, 10, 10, 10, 10),
    @(22, 'bedrooms', 'distractor', 'f', 'stimuli/img_th7xh.png', 82.35897435897436, 65.53846153846153, 73.94871794871796, 39, 8, 8, 8, 8, 8, 8),
    @(23, 'kitchens', 'distractor', 'f', 'stimuli/img_7wquy.png', 50.59375, 30.59375, 40.59375, 32, 2, 2, 2, 2, 2, 2),
    @(24, 'bedrooms', 'distractor', 'f', 'stimuli/img_lgyo6.png', 54.68421052631579, 27.18421052631579, 40.93421052631579, 38, 2, 2, 2, 2, 2, 2),
    @(25, 'living_rooms', 'target', 'j', 'stimuli/img_vgaye.png', 80.33333333333333, 64.57777777777778, 72.45555555555555, 45, 8, 8, 8, 8, 7, 8),
    @(26, 'living_rooms', 'target', 'j', 'stimuli/img_ra2nm.png', 70.75, 50.375, 60.5625, 40, 6, 6, 6, 5, 5, 6),
    @(27, 'living_rooms', 'target', 'j', 'stimuli/img_xr3up.png', 76.24444444444444, 55.88888888888889, 66.06666666666666, 45, 7, 7, 7, 6, 6, 6),
    @(28, 'living_rooms', 'target', 'j', 'stimuli/img_q9lab.png', 53.97560975609756, 32.90243902439025, 43.4390243902439, 41, 3, 3, 3, 3, 4, 3),
    @(29, 'living_rooms', 'target', 'j', 'stimuli/img_pjfx6.png', 32.23404255319149, 26.59574468085106, 29.41489361702127, 47, 2, 2, 2, 2, 2, 3),
    @(30, 'bedrooms', 'distractor', 'f', 'stimuli/img_gka64.png', 19.23809523809524, 20.02380952380953, 19.63095238095238, 42, 1, 1, 1, 1, 1, 2),
    @(31, 'kitchens', 'distractor', 'f', 'stimuli/img_uegbb.png', 78.80952380952381, 61.52380952380953, 70.16666666666667, 42, 8, 8, 8, 8, 8, 8),
    @(32, 'living_rooms', 'target', 'j', 'stimuli/img_c89x3.png', 72.8695652173913, 49.65217391304348, 61.26086956521739, 46, 6, 6, 6, 6, 6, 5),
    @(33, 'living_rooms', 'target', 'j', 'stimuli/img_j856a.png', 38.225, 25.875, 32.05, 40, 2, 2, 2, 3, 3, 2),
    @(34, 'living_rooms', 'target', 'j', 'stimuli/img_wbws6.png', 57.97777777777777, 42.53333333333333, 50.25555555555555, 45, 4, 4, 4, 4, 4, 5),
    @(35, 'living_rooms', 'target', 'j', 'stimuli/img_b21d7.png', 27.75555555555556, 13.86666666666667, 20.81111111111111, 45, 1, 1, 1, 2, 2, 1),
    @(36, 'living_rooms', 'target', 'j', 'stimuli/img_tn8ys.png', 86.70454545454545, 72.40909090909091, 79.55681818181819, 44, 10, 10, 10, 9, 9, 10),
    @(37, 'living_rooms', 'target', 'j', 'stimuli/img_x4bln.png', 76.34042553191489, 59.51063829787234, 67.92553191489361, 47, 7, 7, 7, 7, 7, 7),
    @(38, 'kitchens', 'distractor', 'f', 'stimuli/img_hmmra.png', 54.65853658536585, 34.24390243902439, 44.45121951219512, 41, 3, 3, 3, 4, 4, 3),
    @(39, 'kitchens', 'distractor', 'f', 'stimuli/img_ub9nn.png', 78.77142857142857, 60.37142857142857, 69.57142857142857, 35, 8, 8, 8, 8, 8, 8),
    @(40, 'kitchens', 'distractor', 'f', 'stimuli/img_uhmld.png', 54.19354838709678, 32, 43.09677419354838, 31, 2, 2, 2, 2, 2, 2),
    @(41, 'kitchens', 'distractor', 'f', 'stimuli/img_c30d1.png', 78.875, 60.34375, 69.609375, 32, 8, 8, 8, 8, 8, 8)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 8).Value  = $row[1]   # H category
    $ws.Cells.Item($r, 9).Value  = $row[2]   # I cond_cat
    $ws.Cells.Item($r, 11).Value = $row[3]   # K correct_answer
    $ws.Cells.Item($r, 12).Value = $row[4]   # L stimulus
    $ws.Cells.Item($r, 13).Value = $row[5]   # M conceptual
    $ws.Cells.Item($r, 14).Value = $row[6]   # N perceptual
    $ws.Cells.Item($r, 15).Value = $row[7]   # O typicality
    $ws.Cells.Item($r, 16).Value = $row[8]   # P n
    $ws.Cells.Item($r, 17).Value = $row[9]   # Q p_typicality
    $ws.Cells.Item($r, 18).Value = $row[10]  # R p_conceptual
    $ws.Cells.Item($r, 19).Value = $row[11]  # S p_perceptual
    $ws.Cells.Item($r, 20).Value = $row[12]  # T r_typicality
    $ws.Cells.Item($r, 21).Value = $row[13]  # U r_conceptual
    $ws.Cells.Item($r, 22).Value = $row[14]  # V r_perceptual
}
